$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.25"
$ws.Range("E2").Value = "'-0.13%"
$ws.Range("F2").Value = "'8-1-2023"
$ws.Range("G2").Value = "'1"

$ws.Range("D3").Value = "'27.40"
$ws.Range("E3").Value = "'1.43%"
$ws.Range("F3").Value = "'8-1-2023"
$ws.Range("G3").Value = "'1"

$ws.Range("D4").Value = "'4.707"
$ws.Range("E4").Value = "'0.38%"
$ws.Range("F4").Value = "'8-1-2023"
$ws.Range("G4").Value = "'1"

$ws.Range("D5").Value = "'0.06244"
$ws.Range("E5").Value = "'2.71%"
$ws.Range("F5").Value = "'8-1-2023"
$ws.Range("G5").Value = "'1"

$ws.Range("D6").Value = "'6.721"
$ws.Range("E6").Value = "'0.16%"
$ws.Range("F6").Value = "'8-1-2023"
$ws.Range("G6").Value = "'1"

$ws.Range("D7").Value = "'0.8523"
$ws.Range("E7").Value = "'-1.34%"
$ws.Range("F7").Value = "'8-1-2023"
$ws.Range("G7").Value = "'1"

$ws.Range("D8").Value = "'0.9105"
$ws.Range("E8").Value = "'-0.42%"
$ws.Range("F8").Value = "'8-1-2023"
$ws.Range("G8").Value = "'1"

$ws.Range("D9").Value = "'0.1397"
$ws.Range("E9").Value = "'-0.88%"
$ws.Range("F9").Value = "'8-1-2023"
$ws.Range("G9").Value = "'1"

$ws.Range("D10").Value = "'0.04783"
$ws.Range("E10").Value = "'-2.22%"
$ws.Range("F10").Value = "'8-1-2023"
$ws.Range("G10").Value = "'1"

$ws.Range("D11").Value = "'0.07098"
$ws.Range("E11").Value = "'-0.89%"
$ws.Range("F11").Value = "'8-1-2023"
$ws.Range("G11").Value = "'1"

$ws.Range("D12").Value = "'0.03131"
$ws.Range("E12").Value = "'2.07%"
$ws.Range("F12").Value = "'8-1-2023"
$ws.Range("G12").Value = "'1"

$ws.Range("D13").Value = "'0.09055"
$ws.Range("E13").Value = "'-0.88%"
$ws.Range("F13").Value = "'8-1-2023"
$ws.Range("G13").Value = "'1"

$ws.Range("D14").Value = "'0.001536"
$ws.Range("E14").Value = "'0.17%"
$ws.Range("F14").Value = "'8-1-2023"
$ws.Range("G14").Value = "'1"

$ws.Range("D15").Value = "'0.0006151"
$ws.Range("E15").Value = "'1.27%"
$ws.Range("F15").Value = "'8-1-2023"
$ws.Range("G15").Value = "'1"

$ws.Range("D16").Value = "'0.006143"
$ws.Range("E16").Value = "'-0.84%"
$ws.Range("F16").Value = "'8-1-2023"
$ws.Range("G16").Value = "'1"

$ws.Range("D17").Value = "'3.466"
$ws.Range("E17").Value = "'-0.83%"
$ws.Range("F17").Value = "'8-1-2023"
$ws.Range("G17").Value = "'1"

$ws.Range("D18").Value = "'3.175"
$ws.Range("E18").Value = "'0.27%"
$ws.Range("F18").Value = "'8-1-2023"
$ws.Range("G18").Value = "'1"

$ws.Range("D19").Value = "'2.166"
$ws.Range("E19").Value = "'-1.02%"
$ws.Range("F19").Value = "'8-1-2023"
$ws.Range("G19").Value = "'1"

$ws.Range("F20").Value = "'8-1-2023"
$ws.Range("G20").Value = "'1"

$ws.Range("D21").Value = "'0.1301"
$ws.Range("E21").Value = "'0.16%"
$ws.Range("F21").Value = "'8-1-2023"
$ws.Range("G21").Value = "'1"

$ws.Range("D22").Value = "'4.088"
$ws.Range("E22").Value = "'-0.28%"
$ws.Range("F22").Value = "'8-1-2023"
$ws.Range("G22").Value = "'1"

$ws.Range("D23").Value = "'0.04243"
$ws.Range("E23").Value = "'-0.48%"
$ws.Range("F23").Value = "'8-1-2023"
$ws.Range("G23").Value = "'1"

$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'0.09%"
$ws.Range("F24").Value = "'8-1-2023"
$ws.Range("G24").Value = "'1"

$ws.Range("D25").Value = "'0.004099"
$ws.Range("E25").Value = "'1.04%"
$ws.Range("F25").Value = "'8-1-2023"
$ws.Range("G25").Value = "'1"

$ws.Range("E26").Value = "'0.09%"
$ws.Range("F26").Value = "'8-1-2023"
$ws.Range("G26").Value = "'1"

$ws.Range("E27").Value = "'3.39%"
$ws.Range("F27").Value = "'8-1-2023"
$ws.Range("G27").Value = "'1"

$ws.Range("F28").Value = "'8-1-2023"
$ws.Range("G28").Value = "'1"

$ws.Range("F29").Value = "'8-1-2023"
$ws.Range("G29").Value = "'1"

$ws.Range("F30").Value = "'8-1-2023"
$ws.Range("G30").Value = "'1"

$ws.Range("F31").Value = "'8-1-2023"
$ws.Range("G31").Value = "'1"

$ws.Range("F32").Value = "'8-1-2023"
$ws.Range("G32").Value = "'1"

$ws.Range("F33").Value = "'8-1-2023"
$ws.Range("G33").Value = "'1"

$ws.Range("F34").Value = "'8-1-2023"
$ws.Range("G34").Value = "'1"

$ws.Range("F35").Value = "'8-1-2023"
$ws.Range("G35").Value = "'1"

$ws.Range("F36").Value = "'8-1-2023"
$ws.Range("G36").Value = "'1"

$ws.Range("F37").Value = "'8-1-2023"
$ws.Range("G37").Value = "'1"

$ws.Range("F38").Value = "'8-1-2023"
$ws.Range("G38").Value = "'1"

$ws.Range("F39").Value = "'8-1-2023"
$ws.Range("G39").Value = "'1"

$ws.Range("D40").Value = "'0.03878"
$ws.Range("E40").Value = "'-0.10%"
$ws.Range("F40").Value = "'8-1-2023"
$ws.Range("G40").Value = "'1"

$ws.Range("D41").Value = "'0.1110"
$ws.Range("E41").Value = "'-0.64%"
$ws.Range("F41").Value = "'8-1-2023"
$ws.Range("G41").Value = "'1"

$ws.Range("D42").Value = "'0.004117"
$ws.Range("E42").Value = "'-0.60%"
$ws.Range("F42").Value = "'8-1-2023"
$ws.Range("G42").Value = "'1"

$ws.Range("D43").Value = "'0.002184"
$ws.Range("E43").Value = "'-0.76%"
$ws.Range("F43").Value = "'8-1-2023"
$ws.Range("G43").Value = "'1"

$ws.Range("D44").Value = "'0.01357"
$ws.Range("E44").Value = "'-10.27%"
$ws.Range("F44").Value = "'8-1-2023"
$ws.Range("G44").Value = "'1"

$ws.Range("D45").Value = "'0.00005150"
$ws.Range("E45").Value = "'-2.56%"
$ws.Range("F45").Value = "'8-1-2023"
$ws.Range("G45").Value = "'1"

$ws.Range("E46").Value = "'0.08%"
$ws.Range("F46").Value = "'8-1-2023"
$ws.Range("G46").Value = "'1"

$ws.Range("D47").Value = "'0.03503"
$ws.Range("F47").Value = "'8-1-2023"
$ws.Range("G47").Value = "'1"

$ws.Range("D48").Value = "'0.05565"
$ws.Range("E48").Value = "'-57.87%"
$ws.Range("F48").Value = "'8-1-2023"
$ws.Range("G48").Value = "'1"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.08%"
$ws.Range("F49").Value = "'8-1-2023"
$ws.Range("G49").Value = "'1"

$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.08%"
$ws.Range("F50").Value = "'8-1-2023"
$ws.Range("G50").Value = "'1"

$ws.Range("F51").Value = "'8-1-2023"
$ws.Range("G51").Value = "'1"
